$wb = $excel.ActiveWorkbook

# "Chart" sheet holds the daily GSC export rows (Date, Invalid, Valid).
$ws = $wb.Worksheets.Item("Chart")

# Append the new day's row right after the last existing data row.
# A leading apostrophe forces the date-like text to stay a plain string
# (matching every other row in column A) instead of being auto-converted
# into a serial date value/format.
$ws.Range("A79").Value2 = "'2025-12-22"
$ws.Range("B79").Value2 = 0
$ws.Range("C79").Value2 = 32
